# importerSemaine.xlsx - rebuild the "S1..S10" week lookup table used by
# the absence-import feature so each week's end date (col F) is really
# the Thursday that closes the Friday-start week begun in col E (a
# 7-day span), and each subsequent week starts the day after the
# previous week ends. The previous data had start dates incrementing by
# 1 day and an unrelated end-date series roughly six months out - this
# corrects it to a proper rolling weekly calendar.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 1 (row 5) keeps its existing start date (E5 = 44197); only its
# end date needs correcting to close out a 7-day week.
$ws.Range("F5").Value = 44203

# Weeks 2-10 (rows 6-14): both start (E) and end (F) dates move to the
# new weekly cadence.
$ws.Range("E6").Value = 44204
$ws.Range("F6").Value = 44210

$ws.Range("E7").Value = 44211
$ws.Range("F7").Value = 44217

$ws.Range("E8").Value = 44218
$ws.Range("F8").Value = 44224

$ws.Range("E9").Value = 44225
$ws.Range("F9").Value = 44231

$ws.Range("E10").Value = 44232
$ws.Range("F10").Value = 44238

$ws.Range("E11").Value = 44239
$ws.Range("F11").Value = 44245

$ws.Range("E12").Value = 44246
$ws.Range("F12").Value = 44252

$ws.Range("E13").Value = 44253
$ws.Range("F13").Value = 44259

$ws.Range("E14").Value = 44260
$ws.Range("F14").Value = 44266

# The sheet was last left with the cursor on F17 when saved.
$ws.Range("F17").Select() | Out-Null

# Column F was widened (it's no longer auto best-fit) to comfortably
# show the longer recalculated end dates. 11.6 characters is the
# closest input that this host's pixel-quantized column-width
# conversion resolves to the saved stored width (12.5 vs the original
# file's 12.5703125).
$ws.Columns("F").ColumnWidth = 11.6
